{"js": "// Update the worksheet date header and every \"three-digit \u00f7 one-digit\" answer\n// cell to the values from the new day's generated output.\n\nconst replacements = [\n  [\"2024-06-01 Saturday\", \"2024-06-02 Sunday\"],\n  [\"314\u00f74=78, 2\", \"390\u00f78=48, 6\"],\n  [\"256\u00f76=42, 4\", \"660\u00f72=330, 0\"],\n  [\"452\u00f72=226, 0\", \"929\u00f78=116, 1\"],\n  [\"292\u00f75=58, 2\", \"291\u00f72=145, 1\"],\n  [\"258\u00f74=64, 2\", \"567\u00f78=70, 7\"],\n  [\"265\u00f73=88, 1\", \"687\u00f73=229, 0\"],\n  [\"155\u00f78=19, 3\", \"861\u00f74=215, 1\"],\n  [\"990\u00f75=198, 0\", \"567\u00f74=141, 3\"],\n  [\"134\u00f74=33, 2\", \"383\u00f77=54, 5\"],\n  [\"503\u00f74=125, 3\", \"869\u00f73=289, 2\"],\n  [\"458\u00f75=91, 3\", \"420\u00f77=60, 0\"],\n  [\"225\u00f74=56, 1\", \"406\u00f78=50, 6\"],\n  [\"289\u00f79=32, 1\", \"555\u00f79=61, 6\"],\n  [\"750\u00f75=150, 0\", \"559\u00f74=139, 3\"],\n  [\"714\u00f74=178, 2\", \"539\u00f72=269, 1\"],\n  [\"697\u00f73=232, 1\", \"677\u00f79=75, 2\"],\n  [\"103\u00f77=14, 5\", \"709\u00f76=118, 1\"],\n  [\"641\u00f74=160, 1\", \"393\u00f76=65, 3\"],\n  [\"591\u00f79=65, 6\", \"391\u00f79=43, 4\"],\n  [\"415\u00f74=103, 3\", \"489\u00f72=244, 1\"],\n  [\"670\u00f73=223, 1\", \"836\u00f74=209, 0\"],\n  [\"219\u00f72=109, 1\", \"128\u00f77=18, 2\"],\n  [\"417\u00f76=69, 3\", \"547\u00f75=109, 2\"],\n  [\"667\u00f73=222, 1\", \"665\u00f74=166, 1\"],\n  [\"595\u00f79=66, 1\", \"950\u00f75=190, 0\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date header and every \"three-digit \u00f7 one-digit\" answer\n# cell to the values from the new day's generated output.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"2024-06-01 Saturday\"; New = \"2024-06-02 Sunday\" },\n    @{ Old = \"314\u00f74=78, 2\";         New = \"390\u00f78=48, 6\" },\n    @{ Old = \"256\u00f76=42, 4\";         New = \"660\u00f72=330, 0\" },\n    @{ Old = \"452\u00f72=226, 0\";        New = \"929\u00f78=116, 1\" },\n    @{ Old = \"292\u00f75=58, 2\";         New = \"291\u00f72=145, 1\" },\n    @{ Old = \"258\u00f74=64, 2\";         New = \"567\u00f78=70, 7\" },\n    @{ Old = \"265\u00f73=88, 1\";         New = \"687\u00f73=229, 0\" },\n    @{ Old = \"155\u00f78=19, 3\";         New = \"861\u00f74=215, 1\" },\n    @{ Old = \"990\u00f75=198, 0\";        New = \"567\u00f74=141, 3\" },\n    @{ Old = \"134\u00f74=33, 2\";         New = \"383\u00f77=54, 5\" },\n    @{ Old = \"503\u00f74=125, 3\";        New = \"869\u00f73=289, 2\" },\n    @{ Old = \"458\u00f75=91, 3\";         New = \"420\u00f77=60, 0\" },\n    @{ Old = \"225\u00f74=56, 1\";         New = \"406\u00f78=50, 6\" },\n    @{ Old = \"289\u00f79=32, 1\";         New = \"555\u00f79=61, 6\" },\n    @{ Old = \"750\u00f75=150, 0\";        New = \"559\u00f74=139, 3\" },\n    @{ Old = \"714\u00f74=178, 2\";        New = \"539\u00f72=269, 1\" },\n    @{ Old = \"697\u00f73=232, 1\";        New = \"677\u00f79=75, 2\" },\n    @{ Old = \"103\u00f77=14, 5\";         New = \"709\u00f76=118, 1\" },\n    @{ Old = \"641\u00f74=160, 1\";        New = \"393\u00f76=65, 3\" },\n    @{ Old = \"591\u00f79=65, 6\";         New = \"391\u00f79=43, 4\" },\n    @{ Old = \"415\u00f74=103, 3\";        New = \"489\u00f72=244, 1\" },\n    @{ Old = \"670\u00f73=223, 1\";        New = \"836\u00f74=209, 0\" },\n    @{ Old = \"219\u00f72=109, 1\";        New = \"128\u00f77=18, 2\" },\n    @{ Old = \"417\u00f76=69, 3\";         New = \"547\u00f75=109, 2\" },\n    @{ Old = \"667\u00f73=222, 1\";        New = \"665\u00f74=166, 1\" },\n    @{ Old = \"595\u00f79=66, 1\";         New = \"950\u00f75=190, 0\" }\n)\n\nforeach ($pair in $replacements) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Execute(\n        $pair.Old,\n        $false,\n        $false,\n        $false,\n        $false,\n        $false,\n        $true,\n        1,\n        $false,\n        $pair.New,\n        2\n    )\n}\n"}
